$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers must be forced to
# Text format first, otherwise Excel auto-converts them to numeric
# values (losing precision / the original text formatting).
$ws.Range("D2").Value = '27.850.42'
$ws.Range("E2").Value = '  +1.42%  '
$ws.Range("D3").Value = '1.640.42'
$ws.Range("E3").Value = '  +1.59%  '
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.83'
$ws.Range("E5").Value = '  +0.87%  '
$ws.Range("E6").Value = '  +0.11%  '
$ws.Range("E7").Value = '  -0.20%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.49'
$ws.Range("E8").Value = '  +2.60%  '
$ws.Range("E9").Value = '  +2.15%  '
$ws.Range("E11").Value = '  -1.96%  '
$ws.Range("D12").Value = '1.872.63'
$ws.Range("E12").Value = '  +1.45%  '
$ws.Range("D13").Value = '1.643.11'
$ws.Range("E13").Value = '  +2.73%  '
$ws.Range("E14").Value = '  +1.36%  '
$ws.Range("E15").Value = '  +2.53%  '
$ws.Range("E16").Value = '  +2.20%  '
$ws.Range("D17").Value = '27.858.15'
$ws.Range("E17").Value = '  +1.40%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '231.56'
$ws.Range("E18").Value = '  +1.99%  '
$ws.Range("E19").Value = '  +2.03%  '
$ws.Range("D20").Value = '0.0₃0722'
$ws.Range("E20").Value = '  +0.52%  '
$ws.Range("E21").Value = '  -0.13%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.72'
$ws.Range("E22").Value = '  +8.55%  '
$ws.Range("E23").Value = '  +2.44%  '
$ws.Range("E24").Value = '  +4.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.49'
$ws.Range("E25").Value = '  +1.69%  '
$ws.Range("E26").Value = '  +1.55%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.69'
$ws.Range("E27").Value = '  +1.07%  '
$ws.Range("B28").Value = 'Stellar'
$ws.Range("C28").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.111'
$ws.Range("E28").Value = '  +0.10%  '
$ws.Range("E29").Value = '  -0.18%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.19'
$ws.Range("E30").Value = '  +0.77%  '
$ws.Range("E31").Value = '  +0.46%  '
$ws.Range("E32").Value = '  +1.35%  '
$ws.Range("D33").Value = '1.458.04'
$ws.Range("E33").Value = '  +1.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.11'
$ws.Range("E34").Value = '  +1.97%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.56'
$ws.Range("E35").Value = '  +2.21%  '
$ws.Range("E36").Value = '  -0.64%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.888'
$ws.Range("E37").Value = '  +2.99%  '
$ws.Range("E38").Value = '  +0.15%  '
$ws.Range("E39").Value = '  +0.82%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.919'
$ws.Range("E40").Value = '  -2.62%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '69.22'
$ws.Range("E41").Value = '  +0.21%  '
$ws.Range("E42").Value = '  -0.18%  '
$ws.Range("E43").Value = '  +0.58%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.47'
$ws.Range("E44").Value = '  +0.14%  '
$ws.Range("E45").Value = '  +0.71%  '
$ws.Range("E46").Value = '  -0.08%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.77'
$ws.Range("E47").Value = '  +5.84%  '
$ws.Range("D48").Value = '1.781.66'
$ws.Range("E48").Value = '  +1.25%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '88.40'
$ws.Range("E49").Value = '  +2.70%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.101'
$ws.Range("E50").Value = '  +2.00%  '
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").Value = '0.0₇0997'
$ws.Range("E51").Value = '  -5.27%  '
